$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data block of weekly Cilantro price records occupies rows 122-154
# (columns A-R). This edit inserts one new weekly record by shifting the
# varying columns (D, J, K, L, M, P) of every existing row down by one
# row, moving the former row 154 values into a brand-new row 155, and
# placing a brand-new record's values into row 122 (the first row of the
# block, which becomes the most recent entry after the shift).
#
# Note: this runtime's Range.Value getter is not reliable for reading
# values back out (it yields a property-descriptor placeholder instead
# of the underlying data when stored/reused), so Value2 is used for all
# reads/writes instead, which behaves correctly.

$firstRow = 122
$lastRow = 154
$newRow = 155

# Columns that vary row to row within this block.
$varCols = @("D", "J", "K", "L", "M", "P")

# Columns that stay constant across the whole block; used to populate the
# brand new row 155 (copied straight from the last existing row, 154).
$constCols = @("A", "B", "C", "E", "F", "G", "H", "I", "N", "O", "Q", "R")

# 1) Snapshot the current (pre-edit) values for the varying columns of
#    every row in the block, before any writes happen.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    foreach ($col in $varCols) {
        $rowVals[$col] = $ws.Range("$col$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# 2) Shift rows 123..154 down: row r gets what row (r-1) used to contain.
for ($r = $lastRow; $r -ge ($firstRow + 1); $r--) {
    $prev = $snapshot[$r - 1]
    foreach ($col in $varCols) {
        $ws.Range("$col$r").Value2 = $prev[$col]
    }
}

# 3) Build the brand new row 155, copying constant columns from row 154
#    (identical across the whole block) and the varying values that used
#    to sit in row 154 before the shift.
foreach ($col in $constCols) {
    $ws.Range("$col$newRow").Value2 = $ws.Range("$col$lastRow").Value2
}
$ws.Range("D$newRow").NumberFormat = $ws.Range("D$lastRow").NumberFormat

$last = $snapshot[$lastRow]
foreach ($col in $varCols) {
    $ws.Range("$col$newRow").Value2 = $last[$col]
}

# 4) Write the brand new weekly record's values into row 122 (now the
#    top/most-recent row of the block).
$ws.Range("D$firstRow").Value2 = 44754
$ws.Range("J$firstRow").Value2 = 3200
$ws.Range("K$firstRow").Value2 = 1500
$ws.Range("L$firstRow").Value2 = 2000
$ws.Range("M$firstRow").Value2 = 1750
$ws.Range("P$firstRow").Value2 = 1167
